$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 0.72790054227419165
$ws.Range("C2").Value = 0.93662520785023462
$ws.Range("D2").Value = 0.48922537474695921
$ws.Range("E2").Value = 1.0106203106601002

# Row 3 values
$ws.Range("B3").Value = 0.25390981566173826
$ws.Range("C3").Value = 1.269633190629619
$ws.Range("D3").Value = 0.77385520587044976
$ws.Range("E3").Value = 0.98601100741103531

# Update selection to reflect new active range
$ws.Range("B1:E3").Select()
